$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13, shifting existing rows 13..44 down to 14..45.
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the new record (same Mercado/Categoria metadata,
# new Fecha/Volumen/Precio/Unidad/Origen data).
$ws.Range("A13").Value = 7
$ws.Range("B13").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C13").Value = "Ñuble"
$ws.Range("D13").Value = 44544
$ws.Range("E13").Value = 16
$ws.Range("F13").Value = 100112021
$ws.Range("G13").Value = "Ají"
$ws.Range("H13").Value = "Americana (o)"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 17500
$ws.Range("L13").Value = 18000
$ws.Range("M13").Value = 17750
$ws.Range("N13").Value = "$/caja 15 kilos"
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 1183
$ws.Range("Q13").Value = 15
$ws.Range("R13").Value = "Hortaliza"
